# ChangedTags     Seq.toList    -> Seq.toArray
#
# The underlying change is internal to the exporter that produced this
# deck: the list of "ChangedTags" driving which SharePoint / document-
# management custom XML parts get (re)written was switched from
# `Seq.toList` to `Seq.toArray`. That changes nothing about the *content*
# of the two affected custom XML parts - only the order they are
# (re)materialised into the package, which flips which physical part
# name each one lands in:
#
#   customXml/item2.xml       content-type schema  <->  properties
#   customXml/item3.xml       properties            <->  content-type schema
#
# customXml/itemProps2.xml / itemProps3.xml (the matching datastoreItem
# schemaRefs, keyed off the ds:itemID of whichever part they describe)
# swap the same way, in lock-step with the item they belong to.
#
# Reproduce the intent via the real CustomXMLParts object model: pull the
# two parts out (by namespace, falling back to their well-known
# datastoreItem ids), capture their XML, delete them, then re-add them in
# the opposite order - "properties" first (so it becomes item2.xml), then
# the content-type schema (so it becomes item3.xml) - mirroring the
# Seq.toArray-ordered re-emit the exporter now performs.

try {
    $p = $ppt.ActivePresentation

    $ctNamespace    = "http://schemas.microsoft.com/office/2006/metadata/contentType"
    $propsNamespace = "http://schemas.microsoft.com/office/2006/metadata/properties"

    $ctItemId    = "{DE4876F9-7AE1-498D-B8FE-1E3AD703D2AF}"
    $propsItemId = "{50811A92-D464-4AC4-A396-BA73B10CEEAC}"

    function Find-CustomXmlPart($ns, $id) {
        $part = $null
        try {
            $found = $p.CustomXMLParts.SelectByNamespace($ns)
            if ($found -and $found.Count -ge 1) { $part = $found.Item(1) }
        } catch { $part = $null }

        if (-not $part) {
            try {
                $byId = $p.CustomXMLParts.SelectByID($id)
                if ($byId -and $byId.Count -ge 1) { $part = $byId.Item(1) }
                elseif ($byId) { $part = $byId }
            } catch { }
        }
        return $part
    }

    $ctPart    = Find-CustomXmlPart $ctNamespace    $ctItemId
    $propsPart = Find-CustomXmlPart $propsNamespace $propsItemId

    $ctXml    = $null
    $propsXml = $null
    try { if ($ctPart)    { $ctXml    = $ctPart.XML } }    catch { $ctXml    = $null }
    try { if ($propsPart) { $propsXml = $propsPart.XML } } catch { $propsXml = $null }

    # Remove the existing parts; their text content is already captured above.
    try { if ($ctPart)    { $ctPart.Delete() } }    catch { }
    try { if ($propsPart) { $propsPart.Delete() } } catch { }

    # Re-add "properties" before the content-type schema so the package
    # numbers them in the swapped order (properties -> item2.xml,
    # content-type schema -> item3.xml).
    try { if ($propsXml) { $p.CustomXMLParts.Add($propsXml) | Out-Null } } catch { }
    try { if ($ctXml)    { $p.CustomXMLParts.Add($ctXml) | Out-Null } }    catch { }
} catch {
    # Never fail the whole edit over this metadata-only, best-effort swap.
}
